$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.431.76'
$ws.Range("E2").Value = '  -4.38%  '
$ws.Range("D3").Value = '2.967.14'
$ws.Range("E3").Value = '  -5.52%  '
$ws.Range("D5").Value = '''537.46'
$ws.Range("E5").Value = '  -5.87%  '
$ws.Range("D6").Value = '''149.71'
$ws.Range("E6").Value = '  -8.14%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '''0.567'
$ws.Range("E8").Value = '  -0.98%  '
$ws.Range("D9").Value = '2.972.72'
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").Value = '''0.113'
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("D11").Value = '''6.13'
$ws.Range("E11").Value = '  -6.83%  '
$ws.Range("E12").Value = '  -4.52%  '
$ws.Range("D13").Value = '3.485.85'
$ws.Range("E13").Value = '  -5.68%  '
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = '61.545.02'
$ws.Range("E15").Value = '  -4.25%  '
$ws.Range("D16").Value = '''23.56'
$ws.Range("E16").Value = '  -6.07%  '
$ws.Range("D17").Value = '2.983.38'
$ws.Range("E17").Value = '  -5.20%  '
$ws.Range("E18").Value = '  -4.85%  '
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("D20").Value = '''11.97'
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("D21").Value = '''379.52'
$ws.Range("E21").Value = '  -5.21%  '
$ws.Range("E22").Value = '  -5.92%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D25").Value = '''65.43'
$ws.Range("E25").Value = '  -4.33%  '
$ws.Range("E26").Value = '  -2.95%  '
$ws.Range("D27").Value = '3.094.19'
$ws.Range("E27").Value = '  -5.41%  '
$ws.Range("D28").Value = '''0.188'
$ws.Range("E28").Value = '  -3.52%  '
$ws.Range("D29").Value = '''0.997'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").Value = '0.0₃0934'
$ws.Range("E30").Value = '  -7.14%  '
$ws.Range("D31").Value = '''8.19'
$ws.Range("E31").Value = '  -6.56%  '
$ws.Range("E33").Value = '  -4.97%  '
$ws.Range("D34").Value = '''20.35'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("D35").Value = '''161.02'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").Value = '''4.64'
$ws.Range("E36").Value = '  -3.64%  '
$ws.Range("E37").Value = '  -5.58%  '
$ws.Range("E38").Value = '  -3.39%  '
$ws.Range("E39").Value = '  -5.65%  '
$ws.Range("E40").Value = '  -7.99%  '
$ws.Range("D41").Value = '''37.48'
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("E42").Value = '  -3.89%  '
$ws.Range("D43").Value = '2.404.46'
$ws.Range("E43").Value = '  -9.14%  '
$ws.Range("D44").Value = '''22.01'
$ws.Range("E44").Value = '  -6.78%  '
$ws.Range("D45").Value = '''0.666'
$ws.Range("E45").Value = '  -3.44%  '
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("D47").Value = '''5.08'
$ws.Range("E47").Value = '  -6.54%  '
$ws.Range("D48").Value = '''0.997'
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").Value = '''0.0245'
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("D50").Value = '''0.0950'
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("D51").Value = '''19.61'
$ws.Range("E51").Value = '  -6.56%  '
